$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update headers (renamed, same relative column order)
$ws.Range("A1").Value = "user_id"
$ws.Range("B1").Value = "score"

# Update score (column B) values row by row
$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 8.5
$ws.Range("B4").Value = 6.5
$ws.Range("B5").Value = 9
$ws.Range("B6").Value = 8.5
$ws.Range("B7").Value = 6
$ws.Range("B8").Value = 5.5
$ws.Range("B9").Value = 5
$ws.Range("B10").Value = 7
$ws.Range("B11").ClearContents()
$ws.Range("B12").Value = 8.5
$ws.Range("B13").Value = 7
$ws.Range("B14").Value = 8
$ws.Range("B15").Value = 6.5
$ws.Range("B16").Value = 5
$ws.Range("B17").Value = 8.5
$ws.Range("B18").Value = 4.5
$ws.Range("B19").Value = 6
$ws.Range("B20").Value = 7
$ws.Range("B21").Value = 8
$ws.Range("B22").Value = 9
$ws.Range("B23").Value = 8
$ws.Range("B24").Value = 5
$ws.Range("B25").Value = 6
$ws.Range("B26").Value = 7
$ws.Range("B27").Value = 8
$ws.Range("B28").Value = 10

# Row 28 column A loses its shared formula and becomes a plain value
$ws.Range("A28").Value = 30

# Update the view selection/scroll position to match the edited file
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("A28").Select() | Out-Null
